# Append IUCN status rows for the 5 turtle species to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for new rows 15-19: group, genus, species, common, code, NA
$newRows = @(
    @("bycatch", "Caretta",      "caretta",  "Loggerhead turtle",   "Caretta_caretta_IUCN",      "NA"),
    @("bycatch", "Chelonia",     "mydas",    "Green turtle",        "Chelonia_mydas_IUCN",       "NA"),
    @("bycatch", "Dermochelys",  "coriacea", "Leatherback turtle",  "Dermochelys_coriacea_IUCN", "NA"),
    @("bycatch", "Eretmochelys", "imbricata","Hawksbill turtle",    "Eretmochelys_imbricata_IUCN","NA"),
    @("bycatch", "Lepidochelys", "olivacea", "Olive ridley turtle", "Lepidochelys_olivacea_IUCN","NA")
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Update the active selection to reflect the new last-used cell
$ws.Range("F20").Select()

$wb.Save()
